# Apply the "sleva" (discount) column re-formatting on the KLIENTI sheet:
# values were stored as fractions under a 0.00% format; they're rewritten
# as the equivalent "times 100" plain numbers under a 0.00 number format.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("STROJE")
$ws2 = $wb.Worksheets.Item("KLIENTI")

# Re-style + rewrite column D (discount %) on KLIENTI: D2:D11 move from the
# percentage style to a plain 2-decimal number style, with values scaled by
# 100 so the displayed numbers stay the same.
$rng = $ws2.Range("D2:D11")
$rng.NumberFormat = "0.00"

$ws2.Range("D2").Value = 0
$ws2.Range("D3").Value = 0.05
$ws2.Range("D4").Value = 3.21
$ws2.Range("D5").Value = 0
$ws2.Range("D6").Value = 0.03
$ws2.Range("D7").Value = 4.5999999999999996
$ws2.Range("D8").Value = 58
$ws2.Range("D9").Value = 0.02
$ws2.Range("D10").Value = 0
$ws2.Range("D11").Value = 1.2

# Move the (inactive) STROJE sheet's remembered selection from G24 to D1,
# then restore KLIENTI as the active/visible tab (its own selection, C20,
# is untouched).
$ws1.Range("D1").Select()
$ws2.Activate()
